$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.324023666666666
$ws.Range("H2").Value = 3.972071
$ws.Range("I2").Value = 0.01518042398701374
$ws.Range("J2").Value = 0.01518042398701374
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 153.5290173333333
$ws.Range("N2").Value = 460.587052
$ws.Range("O2").Value = 0.3172206968818489
$ws.Range("P2").Value = 0.317220696881849
$ws.Range("Q2").Value = 203.2760524694102
$ws.Range("R2").Value = 1829.484472224692
$ws.Range("S2").Value = 0.004815544676122433
$ws.Range("T2").Value = 0.004815544676122434
# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.324023666666666
$ws.Range("H3").Value = 3.972071
$ws.Range("I3").Value = 0.01518042398701374
$ws.Range("J3").Value = 0.01518042398701374
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 168.7997026666667
$ws.Range("N3").Value = 506.3991080000001
$ws.Range("O3").Value = 0.3487728915577651
$ws.Range("P3").Value = 0.3487728915577651
$ws.Range("Q3").Value = 223.4948012569631
$ws.Range("R3").Value = 2011.453211312668
$ws.Range("S3").Value = 0.005294520369023639
$ws.Range("T3").Value = 0.005294520369023639
# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.324023666666666
$ws.Range("H4").Value = 3.972071
$ws.Range("I4").Value = 0.01518042398701374
$ws.Range("J4").Value = 0.01518042398701374
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 68.09032333333333
$ws.Range("N4").Value = 204.27097
$ws.Range("O4").Value = 0.1406878008722904
$ws.Range("P4").Value = 0.1406878008722904
$ws.Range("Q4").Value = 90.15319956431887
$ws.Range("R4").Value = 811.37879607887
$ws.Range("S4").Value = 0.002135700467041929
$ws.Range("T4").Value = 0.00213570046704193
# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 1.324023666666666
$ws.Range("H5").Value = 3.972071
$ws.Range("I5").Value = 0.01518042398701374
$ws.Range("J5").Value = 0.01518042398701374
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 93.562673
$ws.Range("N5").Value = 280.688019
$ws.Range("O5").Value = 0.1933186106880956
$ws.Range("P5").Value = 0.1933186106880956
$ws.Range("Q5").Value = 123.8791933685943
$ws.Range("R5").Value = 1114.912740317349
$ws.Range("S5").Value = 0.002934658474825737
$ws.Range("T5").Value = 0.002934658474825737
# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 81.17653533333333
$ws.Range("H6").Value = 243.529606
$ws.Range("I6").Value = 0.9307191821270077
$ws.Range("J6").Value = 0.9307191821270075
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 153.5290173333333
$ws.Range("N6").Value = 460.587052
$ws.Range("O6").Value = 0.3172206968818489
$ws.Range("P6").Value = 0.317220696881849
$ws.Range("Q6").Value = 12462.95370025128
$ws.Range("R6").Value = 112166.5833022615
$ws.Range("S6").Value = 0.2952433875556338
$ws.Range("T6").Value = 0.2952433875556338
# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 81.17653533333333
$ws.Range("H7").Value = 243.529606
$ws.Range("I7").Value = 0.9307191821270077
$ws.Range("J7").Value = 0.9307191821270075
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 168.7997026666667
$ws.Range("N7").Value = 506.3991080000001
$ws.Range("O7").Value = 0.3487728915577651
$ws.Range("P7").Value = 0.3487728915577651
$ws.Range("Q7").Value = 13702.57502777683
$ws.Range("R7").Value = 123323.1752499915
$ws.Range("S7").Value = 0.3246096203787147
$ws.Range("T7").Value = 0.3246096203787146
# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 81.17653533333333
$ws.Range("H8").Value = 243.529606
$ws.Range("I8").Value = 0.9307191821270077
$ws.Range("J8").Value = 0.9307191821270075
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 68.09032333333333
$ws.Range("N8").Value = 204.27097
$ws.Range("O8").Value = 0.1406878008722904
$ws.Range("P8").Value = 0.1406878008722904
$ws.Range("Q8").Value = 5527.336537926424
$ws.Range("R8").Value = 49746.02884133782
$ws.Range("S8").Value = 0.1309408349631054
$ws.Range("T8").Value = 0.1309408349631054
# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 81.17653533333333
$ws.Range("H9").Value = 243.529606
$ws.Range("I9").Value = 0.9307191821270077
$ws.Range("J9").Value = 0.9307191821270075
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 93.562673
$ws.Range("N9").Value = 280.688019
$ws.Range("O9").Value = 0.1933186106880956
$ws.Range("P9").Value = 0.1933186106880956
$ws.Range("Q9").Value = 7595.093630665613
$ws.Range("R9").Value = 68355.84267599051
$ws.Range("S9").Value = 0.1799253392295538
$ws.Range("T9").Value = 0.1799253392295538
# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 1.192675
$ws.Range("H10").Value = 3.578025
$ws.Range("I10").Value = 0.0136744626508778
$ws.Range("J10").Value = 0.0136744626508778
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 153.5290173333333
$ws.Range("N10").Value = 460.587052
$ws.Range("O10").Value = 0.3172206968818489
$ws.Range("P10").Value = 0.317220696881849
$ws.Range("Q10").Value = 183.1102207480333
$ws.Range("R10").Value = 1647.9919867323
$ws.Range("S10").Value = 0.004337822571596272
$ws.Range("T10").Value = 0.004337822571596272
# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 1.192675
$ws.Range("H11").Value = 3.578025
$ws.Range("I11").Value = 0.0136744626508778
$ws.Range("J11").Value = 0.0136744626508778
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 168.7997026666667
$ws.Range("N11").Value = 506.3991080000001
$ws.Range("O11").Value = 0.3487728915577651
$ws.Range("P11").Value = 0.3487728915577651
$ws.Range("Q11").Value = 201.3231853779667
$ws.Range("R11").Value = 1811.9086684017
$ws.Range("S11").Value = 0.004769281879245313
$ws.Range("T11").Value = 0.004769281879245312
# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 1.192675
$ws.Range("H12").Value = 3.578025
$ws.Range("I12").Value = 0.0136744626508778
$ws.Range("J12").Value = 0.0136744626508778
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 68.09032333333333
$ws.Range("N12").Value = 204.27097
$ws.Range("O12").Value = 0.1406878008722904
$ws.Range("P12").Value = 0.1406878008722904
$ws.Range("Q12").Value = 81.20962638158333
$ws.Range("R12").Value = 730.88663743425
$ws.Range("S12").Value = 0.001923830078462268
$ws.Range("T12").Value = 0.001923830078462269
# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 1.192675
$ws.Range("H13").Value = 3.578025
$ws.Range("I13").Value = 0.0136744626508778
$ws.Range("J13").Value = 0.0136744626508778
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 93.562673
$ws.Range("N13").Value = 280.688019
$ws.Range("O13").Value = 0.1933186106880956
$ws.Range("P13").Value = 0.1933186106880956
$ws.Range("Q13").Value = 111.589861020275
$ws.Range("R13").Value = 1004.308749182475
$ws.Range("S13").Value = 0.00264352812157395
$ws.Range("T13").Value = 0.00264352812157395
# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 3.525915333333334
$ws.Range("H14").Value = 10.577746
$ws.Range("I14").Value = 0.04042593123510095
$ws.Range("J14").Value = 0.04042593123510094
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 153.5290173333333
$ws.Range("N14").Value = 460.587052
$ws.Range("O14").Value = 0.3172206968818489
$ws.Range("P14").Value = 0.317220696881849
$ws.Range("Q14").Value = 541.3303163271991
$ws.Range("R14").Value = 4871.972846944792
$ws.Range("S14").Value = 0.01282394207849643
$ws.Range("T14").Value = 0.01282394207849643
# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 3.525915333333334
$ws.Range("H15").Value = 10.577746
$ws.Range("I15").Value = 0.04042593123510095
$ws.Range("J15").Value = 0.04042593123510094
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 168.7997026666667
$ws.Range("N15").Value = 506.3991080000001
$ws.Range("O15").Value = 0.3487728915577651
$ws.Range("P15").Value = 0.3487728915577651
$ws.Range("Q15").Value = 595.1734598945077
$ws.Range("R15").Value = 5356.561139050569
$ws.Range("S15").Value = 0.01409946893078153
$ws.Range("T15").Value = 0.01409946893078153
# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 3.525915333333334
$ws.Range("H16").Value = 10.577746
$ws.Range("I16").Value = 0.04042593123510095
$ws.Range("J16").Value = 0.04042593123510094
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 68.09032333333333
$ws.Range("N16").Value = 204.27097
$ws.Range("O16").Value = 0.1406878008722904
$ws.Range("P16").Value = 0.1406878008722904
$ws.Range("Q16").Value = 240.0807150926245
$ws.Range("R16").Value = 2160.72643583362
$ws.Range("S16").Value = 0.005687435363680786
$ws.Range("T16").Value = 0.005687435363680786
# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 3.525915333333334
$ws.Range("H17").Value = 10.577746
$ws.Range("I17").Value = 0.04042593123510095
$ws.Range("J17").Value = 0.04042593123510094
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 93.562673
$ws.Range("N17").Value = 280.688019
$ws.Range("O17").Value = 0.1933186106880956
$ws.Range("P17").Value = 0.1933186106880956
$ws.Range("Q17").Value = 329.8940633583527
$ws.Range("R17").Value = 2969.046570225174
$ws.Range("S17").Value = 0.007815084862142207
$ws.Range("T17").Value = 0.007815084862142205
